$d = $word.ActiveDocument

# Locate the target paragraph: the one that starts with "多云转小雨"
$targetIndex = -1
$emptyIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("多云转小雨")) {
        $targetIndex = $i
        $emptyIndex = $i + 1
        break
    }
}

if ($targetIndex -eq -1) {
    throw "could not locate target paragraph"
}

$pContent = $d.Paragraphs.Item($targetIndex)
$pEmpty = $d.Paragraphs.Item($emptyIndex)

# Range spanning both the content paragraph and the following empty paragraph,
# so the replacement can drop the <w:pPr>/rFonts-hint on both and remove the
# stray <w:proofErr/> markers in one shot.
$r = $d.Range($pContent.Range.Start, $pEmpty.Range.End)

$run1 = '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>多云转小雨，</w:t></w:r>'
$run2 = '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>今</w:t></w:r>'
$run3 = '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>天学习了分支管理，创建了一个分支dev。</w:t></w:r>'
$run4 = '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>使用git创建分支简单又快速。</w:t></w:r>'

# A lone trailing empty <w:p/> at the very end of a replaced range collapses
# away (no-op) because there is nothing after it to separate, so seed the
# second (to-become-empty) paragraph with a one-character placeholder run and
# strip that character off afterwards - that reliably keeps it as a distinct
# paragraph with a bare <w:p/>.
$placeholder = '<w:r><w:t>Z</w:t></w:r>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' + $run1 + $run2 + $run3 + $run4 + '</w:p>' +
  '<w:p>' + $placeholder + '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml) | Out-Null

# Strip the placeholder character back out, leaving a truly empty paragraph.
$newEmpty = $d.Paragraphs.Item($targetIndex + 1)
$er = $newEmpty.Range
$er.MoveEnd(1, -1) | Out-Null
$er.Text = ""
